# Apply the "include diff_scale feature for all in AR estimate" edit.
# Updates numeric values in rows 2-6 (columns E, F, G, M, N, O) and
# appends a new row 7 with the "Disg"/"Var" feature row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 197.7
$ws.Range("F2").Value = 6011.89
$ws.Range("G2").Value = 1959.51
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.82

# --- Row 3 updates ---
$ws.Range("E3").Value = 75.01000000000001
$ws.Range("F3").Value = 3278.5
$ws.Range("G3").Value = 2632.47
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 1.21

# --- Row 4 updates ---
$ws.Range("E4").Value = 4841.28
$ws.Range("F4").Value = 166.01
$ws.Range("G4").Value = 1953.53
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 1.17

# --- Row 5 updates ---
$ws.Range("E5").Value = 45.21
$ws.Range("F5").Value = 556.33
$ws.Range("G5").Value = 185448.25
$ws.Range("M5").Value = 1703.7
$ws.Range("N5").Value = 10.94
$ws.Range("O5").Value = 60.61

# --- Row 6 updates ---
$ws.Range("E6").Value = 2.28
$ws.Range("F6").Value = 18.79
$ws.Range("G6").Value = 1.36
$ws.Range("M6").Value = 9150.26
$ws.Range("N6").Value = 332.08
$ws.Range("O6").Value = 14172.03

# --- New row 7 ---
$ws.Range("A7").Value = "Disg"
$ws.Range("B7").Value = "Var"
$ws.Range("E7").Value = 90.12
$ws.Range("F7").Value = 0.5600000000000001
$ws.Range("G7").Value = 0.1
$ws.Range("H7").Value = 0.5
$ws.Range("I7").Value = 0.5
$ws.Range("J7").Value = 0.1
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.1
$ws.Range("M7").Value = 3297.11
$ws.Range("N7").Value = 333.49
$ws.Range("O7").Value = 14198.59
$ws.Range("P7").Value = 0.5
$ws.Range("Q7").Value = 0.5
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.8
$ws.Range("T7").Value = 0.1
